$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha, as serial date), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg) for each row.
# Row 7 is unchanged and therefore omitted.
$updates = @{
    2  = @{ D = 44406; J = 400; K = 20000; L = 22000; M = 20850; P = 1390 }
    3  = @{ D = 44391; J = 160; K = 20000; L = 20000; M = 20000; P = 1333 }
    4  = @{ D = 44365; J = 580; K = 20000; L = 22000; M = 21103; P = 1407 }
    5  = @{ D = 44483; J = 220; K = 18000; L = 20000; M = 18909; P = 1261 }
    6  = @{ D = 44435; J = 140; K = 21000; L = 23000; M = 21714; P = 1448 }
    8  = @{ D = 44398; J = 130; K = 20000; L = 20000; M = 20000; P = 1333 }
    9  = @{ D = 44396; J = 130; K = 22000; L = 22000; M = 22000; P = 1467 }
    10 = @{ D = 44399; J = 150; K = 22000; L = 22000; M = 22000; P = 1467 }
    11 = @{ D = 44453; J = 280; K = 20000; L = 22000; M = 21286; P = 1419 }
    12 = @{ D = 44476; J = 220; K = 20000; L = 22000; M = 20909; P = 1394 }
    13 = @{ D = 44449; J = 220; K = 22000; L = 24000; M = 23091; P = 1539 }
    14 = @{ D = 44392; J = 220; K = 23000; L = 23000; M = 23000; P = 1533 }
    15 = @{ D = 44400; J = 130; K = 24000; L = 24000; M = 24000; P = 1600 }
}

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $ws.Cells.Item($row, 4).Value = $epoch.AddDays([double]$vals.D)
    $ws.Cells.Item($row, 10).Value = $vals.J
    $ws.Cells.Item($row, 11).Value = $vals.K
    $ws.Cells.Item($row, 12).Value = $vals.L
    $ws.Cells.Item($row, 13).Value = $vals.M
    $ws.Cells.Item($row, 16).Value = $vals.P
}
